$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$cell.NumberFormat = "@"
$cell.Value = '65.146.87'
$ws.Range('E2').Value = '  -1.32%  '
$cell = $ws.Range('D3')
$cell.NumberFormat = "@"
$cell.Value = '2.939.05'
$ws.Range('E3').Value = '  -2.67%  '
$ws.Range('E4').Value = '  -0.10%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = "@"
$cell.Value = '567.25'
$ws.Range('E5').Value = '  -3.14%  '
$cell = $ws.Range('D6')
$cell.NumberFormat = "@"
$cell.Value = '158.67'
$ws.Range('E6').Value = '  +2.11%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  -0.23%  '
$cell = $ws.Range('D9')
$cell.NumberFormat = "@"
$cell.Value = '2.933.98'
$ws.Range('E9').Value = '  -2.71%  '
$cell = $ws.Range('D10')
$cell.NumberFormat = "@"
$cell.Value = '6.71'
$ws.Range('E10').Value = '  -4.18%  '
$ws.Range('E11').Value = '  -3.48%  '
$ws.Range('E12').Value = '  +0.50%  '
$ws.Range('E13').Value = '  -0.10%  '
$cell = $ws.Range('D14')
$cell.NumberFormat = "@"
$cell.Value = '34.29'
$ws.Range('E14').Value = '  -0.40%  '
$ws.Range('E15').Value = '  -0.89%  '
$cell = $ws.Range('D16')
$cell.NumberFormat = "@"
$cell.Value = '65.152.25'
$ws.Range('E16').Value = '  -1.31%  '
$cell = $ws.Range('D17')
$cell.NumberFormat = "@"
$cell.Value = '3.428.80'
$ws.Range('E17').Value = '  -2.56%  '
$ws.Range('E18').Value = '  -0.85%  '
$cell = $ws.Range('D19')
$cell.NumberFormat = "@"
$cell.Value = '2.940.84'
$ws.Range('E19').Value = '  -2.57%  '
$cell = $ws.Range('D20')
$cell.NumberFormat = "@"
$cell.Value = '14.68'
$ws.Range('E20').Value = '  +6.11%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = "@"
$cell.Value = '445.20'
$ws.Range('E21').Value = '  -3.80%  '
$cell = $ws.Range('D22')
$cell.NumberFormat = "@"
$cell.Value = '0.685'
$ws.Range('E22').Value = '  -0.36%  '
$ws.Range('E23').Value = '  -2.29%  '
$cell = $ws.Range('D24')
$cell.NumberFormat = "@"
$cell.Value = '82.00'
$ws.Range('E24').Value = '  -0.16%  '
$ws.Range('E25').Value = '  -2.77%  '
$cell = $ws.Range('D26')
$cell.NumberFormat = "@"
$cell.Value = '12.06'
$cell = $ws.Range('D27')
$cell.NumberFormat = "@"
$cell.Value = '10.03'
$ws.Range('E27').Value = '  -7.99%  '
$ws.Range('E28').Value = '  +0.01%  '
$cell = $ws.Range('D29')
$cell.NumberFormat = "@"
$cell.Value = '8.00'
$ws.Range('E29').Value = '  +0.43%  '
$cell = $ws.Range('D30')
$cell.NumberFormat = "@"
$cell.Value = '2.38'
$ws.Range('E30').Value = '  -1.71%  '
$ws.Range('E31').Value = '  -1.92%  '
$ws.Range('E32').Value = '  -2.81%  '
$cell = $ws.Range('D33')
$cell.NumberFormat = "@"
$cell.Value = '27.07'
$ws.Range('E33').Value = '  +0.07%  '
$ws.Range('E34').Value = '  -1.55%  '
$cell = $ws.Range('D35')
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$ws.Range('E35').Value = '  -0.02%  '
$cell = $ws.Range('D36')
$cell.NumberFormat = "@"
$cell.Value = '0.972'
$ws.Range('E36').Value = '  -2.33%  '
$cell = $ws.Range('D37')
$cell.NumberFormat = "@"
$cell.Value = '5.71'
$ws.Range('E37').Value = '  -1.44%  '
$cell = $ws.Range('D38')
$cell.NumberFormat = "@"
$cell.Value = '49.59'
$ws.Range('E38').Value = '  +0.28%  '
$cell = $ws.Range('D39')
$cell.NumberFormat = "@"
$cell.Value = '44.21'
$ws.Range('E39').Value = '  -1.85%  '
$ws.Range('E40').Value = '  -10.09%  '
$ws.Range('E41').Value = '  -2.32%  '
$cell = $ws.Range('D42')
$cell.NumberFormat = "@"
$cell.Value = '2.83'
$ws.Range('E42').Value = '  -6.61%  '
$ws.Range('E43').Value = '  -1.86%  '
$ws.Range('E44').Value = '  -0.65%  '
$cell = $ws.Range('D45')
$cell.NumberFormat = "@"
$cell.Value = '384.32'
$ws.Range('E45').Value = '  -1.69%  '
$cell = $ws.Range('D46')
$cell.NumberFormat = "@"
$cell.Value = '0.0351'
$ws.Range('E46').Value = '  -1.19%  '
$cell = $ws.Range('D47')
$cell.NumberFormat = "@"
$cell.Value = '2.702.63'
$ws.Range('E47').Value = '  -3.57%  '
$cell = $ws.Range('D48')
$cell.NumberFormat = "@"
$cell.Value = '132.97'
$ws.Range('E48').Value = '  -1.46%  '
$ws.Range('E50').Value = '  +4.38%  '
$ws.Range('E51').Value = '  -0.54%  '
